# Commit: "Grid participants + grid comments"
# The OOXML diff appends four new paragraphs right after the paragraph
# containing "Ejemplo Tuneado Mostrar" (the last paragraph in the body,
# just before the final sectPr):
#   1. an empty paragraph holding only a manual page break
#   2. "Que es Firebase" (indented, with a lastRenderedPageBreak marker)
#   3. "Es una plataforma de desarrollo móvil "
#   4. an English blurb about Firebase, styled in Arial/gray/shaded text
#
# We rebuild the exact OOXML for the four paragraphs (so indentation,
# language tags, proofErr spans, and run formatting match the target
# precisely) and inject it with Range.InsertXML.

$d = $word.ActiveDocument

# Sanity-check: confirm the expected anchor paragraph is present before
# mutating anything (fail fast rather than silently inserting in the
# wrong spot if the document doesn't look as expected).
$finder = $d.Content
$found = $finder.Find.Execute("Ejemplo Tuneado Mostrar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor paragraph 'Ejemplo Tuneado Mostrar' not found; refusing to guess an insertion point."
}

# Use a *fresh* Range collapsed to the very end of the document body
# (right before sectPr) as the insertion point. Reusing the Range object
# returned by Find (even after navigating to its Paragraphs(1)) can make
# InsertXML clobber the preceding paragraph mark, so we deliberately grab
# a brand-new Content range here instead.
$insertionPoint = $d.Content
$insertionPoint.Collapse(0)

$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Que es </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Firebase</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">Es una plataforma de desarrollo móvil </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="6C6C6C"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Firebase provides tools for tracking analytics, reporting and fixing app crashes, creating marketing and product experiment.</w:t></w:r></w:p>'

[void]$insertionPoint.InsertXML($newParagraphsXml)
